$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the "interesting observation" paragraph (Discussion section).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(32)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End

$newPara1 = "An interesting observation about the sex of the passengers. The rate of the survived men in the first and second class is more than 90%, though only 40% percent of women in the first class survived and even less in the second. This shows that women and children were not the first to step on the lifeboats, but that a place there was gained by strength and quickness. As for the third class, some survivors remember that a lot of 3rd class passengers were blocked from getting to the boats. This fact gets confirmed by the dataset which shows such a big difference between the survivors of the 3rd class and the two first."

# Replace everything but the trailing paragraph mark.
$textRng1 = $d.Range($p1Start, $p1End - 1)
$textRng1.Text = $newPara1

# Superscript the two "rd" in "3rd" (ordinal suffix).
$rd1 = $d.Range($p1Start + 431, $p1Start + 433)
$rd1.Font.Superscript = $true

$rd2 = $d.Range($p1Start + 595, $p1Start + 597)
$rd2.Font.Superscript = $true

# ---------------------------------------------------------------------------
# 2. Edit the Conclusion paragraph ("I know that ...").
# ---------------------------------------------------------------------------
$rng = $d.Paragraphs.Item(35).Range
$rng.Find.Execute("really sorry", $false, $false, $false, $false, $false, $true, 0, $false, "sorry", 2)

$rng = $d.Paragraphs.Item(35).Range
$rng.Find.Execute("But this", $false, $false, $false, $false, $false, $true, 0, $false, "Nevertheless, this", 2)

$rng = $d.Paragraphs.Item(35).Range
$rng.Find.Execute("point of view this", $false, $false, $false, $false, $false, $true, 0, $false, "point of view, this", 2)

# ---------------------------------------------------------------------------
# 3. Relocate the "_GoBack" bookmark from its own paragraph into the middle
#    of the word "great" ("gre" | "at") inside the Conclusion paragraph.
# ---------------------------------------------------------------------------
$p35Start = $d.Paragraphs.Item(35).Range.Start
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$bmPos = $p35Start + 290
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------------
# 4. Turn the (now empty, bookmark-less) paragraph right after the
#    Conclusion paragraph into a new closing sentence.
# ---------------------------------------------------------------------------
$p36 = $d.Paragraphs.Item(36)
$p36.Range.Text = "As well as the analysis of the dataset gives us the possibility to understand the scenario of what has happened onboard."
$p36.Format.SpaceAfter = 0
$p36.Format.LineSpacingRule = 1

# ---------------------------------------------------------------------------
# 5. Add one more empty paragraph right after it (keeping the remaining
#    trailing empty paragraphs intact). Splitting with a bare carriage
#    return keeps the new paragraph free of any placeholder run.
# ---------------------------------------------------------------------------
$p36 = $d.Paragraphs.Item(36)
$splitPos = $p36.Range.End
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.Text = [char]13

Write-Output "done"
